$d = $word.ActiveDocument

function Replace-ExactText($findText, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Assign directly to the matched Range's Text so the literal
        # characters (straight apostrophes, hyphens, commas, etc.) are
        # written verbatim instead of going through Find/Replace's
        # autocorrect-aware substitution (which turns ' into a curly quote).
        $r.Text = $newText
    }
}

# 1) "(1) Mineral salts only." -> "(1) Mineral'salts only."
Replace-ExactText "(1) Mineral salts only." "(1) Mineral'salts only."

# 2) "(3) Carbon bioxide and mineral salts only." -> "(3) - Carbon bioxide and mineral salts only."
Replace-ExactText "(3) Carbon bioxide and mineral salts only." "(3) - Carbon bioxide and mineral salts only."

# 3) "13. Jody poured 200 mi of water into container P ..." -> "73. Jody poured 200 mi of water Info container P ..."
Replace-ExactText "13. Jody poured 200 mi of water into container P as shown in the diagram below. She" "73. Jody poured 200 mi of water Info container P as shown in the diagram below. She"

# 4) Remove the paragraph containing the inline drawing (image) entirely.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.Delete()
    }
}

# 5) "(1). She was frying to find out if water nas definite mass. ." -> "(1), She was frying io find out if wafer nas definite mass. ,"
Replace-ExactText "(1). She was frying to find out if water nas definite mass. ." "(1), She was frying io find out if wafer nas definite mass. ,"
